$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.693.70"
$ws.Range("E2").Value = "  +3.15%  "
$ws.Range("D3").Value = "3.125.65"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.00"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.92%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.121.60"
$ws.Range("E8").Value = "  +1.93%  "
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("E10").Value = "  +15.66%  "
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("E13").Value = "  +5.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.06"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.63%  "
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "3.643.70"
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "63.586.18"
$ws.Range("E18").Value = "  +3.08%  "
$ws.Range("D19").Value = "3.126.80"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.26"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.37"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.732"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.30"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +10.19%  "
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.17"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("E33").Value = "  -4.03%  "
$ws.Range("D34").Value = "0.0₃0877"
$ws.Range("E34").Value = "  +11.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.38"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +8.85%  "
$ws.Range("E36").Value = "  +2.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.43"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +16.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.13"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.99"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "454.48"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.13%  "
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("D43").Value = "2.906.45"
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("E45").Value = "  +2.30%  "
$ws.Range("E46").Value = "  +3.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.69"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.66"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.27%  "
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("E51").Value = "  +1.92%  "
